# Generate Report for Handoff
# Update the "Latest Handoff Datetime" (column D, row 5 - the
# 29ac8fbb-bded-4c1f-8da8-166cc0a86eaf entry) on both the zh-cn and
# de-de localization-status worksheets to reflect a new handoff.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D5").Value = "2016-03-10 04:45:51"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D5").Value = "2016-03-10 04:45:55"
